$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 12538.154  # ALC!H70: 13499.75 -> 12538.154
$ws.Cells.Item(70, 9).Value = 1333  # ALC!I70: 1500 -> 1333
$ws.Cells.Item(70, 11).Value = 3999  # ALC!K70: 4500 -> 3999
$ws.Cells.Item(70, 13).Value = -3729  # ALC!M70: -4230 -> -3729

$ws.Cells.Item(73, 8).Value = 12538.154  # ALC!H73: 13499.75 -> 12538.154
$ws.Cells.Item(73, 9).Value = 1333  # ALC!I73: 1500 -> 1333
$ws.Cells.Item(73, 11).Value = 3999  # ALC!K73: 4500 -> 3999
$ws.Cells.Item(73, 13).Value = -3063  # ALC!M73: -3564 -> -3063

$ws.Cells.Item(86, 8).Value = 1619.2  # ALC!H86: 1619.4 -> 1619.2
$ws.Cells.Item(86, 9).Value = 1399.5  # ALC!I86: 1400 -> 1399.5
$ws.Cells.Item(86, 11).Value = 1399.5  # ALC!K86: 1400 -> 1399.5
$ws.Cells.Item(86, 13).Value = -276.5  # ALC!M86: -277 -> -276.5

$ws.Cells.Item(89, 8).Value = 1619.2  # ALC!H89: 1619.4 -> 1619.2
$ws.Cells.Item(89, 9).Value = 1399.5  # ALC!I89: 1400 -> 1399.5
$ws.Cells.Item(89, 11).Value = 6997.5  # ALC!K89: 7000 -> 6997.5
$ws.Cells.Item(89, 13).Value = -1381.5  # ALC!M89: -1384 -> -1381.5

$ws.Cells.Item(106, 8).Value = 2580  # ALC!H106: 1972.2222 -> 2580
$ws.Cells.Item(106, 9).Value = 2580  # ALC!I106: 1972.2222 -> 2580
$ws.Cells.Item(106, 11).Value = 2580  # ALC!K106: 1972.2222 -> 2580
$ws.Cells.Item(106, 13).Value = -1949  # ALC!M106: -1341.2222 -> -1949

$ws.Cells.Item(130, 8).Value = 40000  # ALC!H130: 39593.332 -> 40000
$ws.Cells.Item(130, 10).Value = 40000  # ALC!J130: 39593.332 -> 40000
$ws.Cells.Item(130, 12).Value = 40000  # ALC!L130: 39593.332 -> 40000
$ws.Cells.Item(130, 14).Value = -50040  # ALC!N130: -49633.332 -> -50040

$ws.Cells.Item(132, 8).Value = 1434.1333  # ALC!H132: 1491.2 -> 1434.1333
$ws.Cells.Item(132, 9).Value = 1444.4445  # ALC!I132: 1516.6666 -> 1444.4445
$ws.Cells.Item(132, 10).Value = 1418.6666  # ALC!J132: 1453 -> 1418.6666
$ws.Cells.Item(132, 11).Value = 4333.333500000001  # ALC!K132: 4549.9998 -> 4333.333500000001
$ws.Cells.Item(132, 12).Value = 4255.9998  # ALC!L132: 4359 -> 4255.9998
$ws.Cells.Item(132, 13).Value = -1803.333500000001  # ALC!M132: -2019.9998 -> -1803.333500000001
$ws.Cells.Item(132, 14).Value = -9315.9998  # ALC!N132: -9419 -> -9315.9998

$ws.Cells.Item(137, 8).Value = 1851.35  # ALC!H137: 1557.9032 -> 1851.35
$ws.Cells.Item(137, 9).Value = 1531  # ALC!I137: 1330.3043 -> 1531
$ws.Cells.Item(137, 10).Value = 3666.6667  # ALC!J137: 2212.25 -> 3666.6667
$ws.Cells.Item(137, 11).Value = 4593  # ALC!K137: 3990.9129 -> 4593
$ws.Cells.Item(137, 12).Value = 11000.0001  # ALC!L137: 6636.75 -> 11000.0001
$ws.Cells.Item(137, 13).Value = -2043  # ALC!M137: -1440.9129 -> -2043
$ws.Cells.Item(137, 14).Value = -16100.0001  # ALC!N137: -11736.75 -> -16100.0001

$ws.Cells.Item(138, 8).Value = 1924.2461  # ALC!H138: 2121.7378 -> 1924.2461
$ws.Cells.Item(138, 9).Value = 1703.3715  # ALC!I138: 1961.6111 -> 1703.3715
$ws.Cells.Item(138, 10).Value = 2181.9333  # ALC!J138: 2352.32 -> 2181.9333
$ws.Cells.Item(138, 11).Value = 5110.1145  # ALC!K138: 5884.8333 -> 5110.1145
$ws.Cells.Item(138, 12).Value = 6545.7999  # ALC!L138: 7056.960000000001 -> 6545.7999
$ws.Cells.Item(138, 13).Value = 29.88550000000032  # ALC!M138: -744.8333000000002 -> 29.88550000000032
$ws.Cells.Item(138, 14).Value = -16825.7999  # ALC!N138: -17336.96 -> -16825.7999

$ws.Cells.Item(140, 8).Value = 64032.5  # ALC!H140: 65283.684 -> 64032.5
$ws.Cells.Item(140, 10).Value = 64032.5  # ALC!J140: 65283.684 -> 64032.5
$ws.Cells.Item(140, 12).Value = 64032.5  # ALC!L140: 65283.684 -> 64032.5
$ws.Cells.Item(140, 14).Value = -74392.5  # ALC!N140: -75643.68400000001 -> -74392.5

$ws.Cells.Item(141, 8).Value = 5471.6665  # ALC!H141: 7610.25 -> 5471.6665
$ws.Cells.Item(141, 9).Value = 1874.8  # ALC!I141: 2328.3333 -> 1874.8
$ws.Cells.Item(141, 11).Value = 5624.4  # ALC!K141: 6984.999899999999 -> 5624.4
$ws.Cells.Item(141, 13).Value = -444.3999999999996  # ALC!M141: -1804.999899999999 -> -444.3999999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4593.4106  # ARM!H32: 4816.6416 -> 4593.4106
$ws.Cells.Item(32, 9).Value = 3291.5208  # ARM!I32: 3408.5652 -> 3291.5208
$ws.Cells.Item(32, 10).Value = 12404.75  # ARM!J32: 14069.714 -> 12404.75
$ws.Cells.Item(32, 11).Value = 3291.5208  # ARM!K32: 3408.5652 -> 3291.5208
$ws.Cells.Item(32, 12).Value = 12404.75  # ARM!L32: 14069.714 -> 12404.75
$ws.Cells.Item(32, 13).Value = -3004.5208  # ARM!M32: -3121.5652 -> -3004.5208
$ws.Cells.Item(32, 14).Value = -12978.75  # ARM!N32: -14643.714 -> -12978.75

$ws.Cells.Item(61, 8).Value = 1616.3478  # ARM!H61: 1741.238 -> 1616.3478
$ws.Cells.Item(61, 9).Value = 1466.3529  # ARM!I61: 1621.2 -> 1466.3529
$ws.Cells.Item(61, 11).Value = 1466.3529  # ARM!K61: 1621.2 -> 1466.3529
$ws.Cells.Item(61, 13).Value = -1254.3529  # ARM!M61: -1409.2 -> -1254.3529

$ws.Cells.Item(132, 8).Value = 1562.375  # ARM!H132: 1446.5 -> 1562.375
$ws.Cells.Item(132, 9).Value = 1583.2  # ARM!I132: 1446.5 -> 1583.2
$ws.Cells.Item(132, 10).Value = 1250  # ARM!J132: 0 -> 1250
$ws.Cells.Item(132, 11).Value = 4749.6  # ARM!K132: 4339.5 -> 4749.6
$ws.Cells.Item(132, 12).Value = 3750  # ARM!L132: 0 -> 3750
$ws.Cells.Item(132, 13).Value = -2219.6  # ARM!M132: -1809.5 -> -2219.6
$ws.Cells.Item(132, 14).Value = -8810  # ARM!N132: None -> -8810

$ws.Cells.Item(134, 8).Value = 54995  # ARM!H134: 54945 -> 54995
$ws.Cells.Item(134, 10).Value = 54995  # ARM!J134: 54945 -> 54995
$ws.Cells.Item(134, 12).Value = 54995  # ARM!L134: 54945 -> 54995
$ws.Cells.Item(134, 14).Value = -65135  # ARM!N134: -65085 -> -65135

$ws.Cells.Item(136, 8).Value = 1616.3478  # ARM!H136: 1741.238 -> 1616.3478
$ws.Cells.Item(136, 9).Value = 1466.3529  # ARM!I136: 1621.2 -> 1466.3529
$ws.Cells.Item(136, 11).Value = 4399.0587  # ARM!K136: 4863.6 -> 4399.0587
$ws.Cells.Item(136, 13).Value = -1849.0587  # ARM!M136: -2313.6 -> -1849.0587

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1052.8334  # BSM!H99: 1053 -> 1052.8334
$ws.Cells.Item(99, 10).Value = 1052.8334  # BSM!J99: 1053 -> 1052.8334
$ws.Cells.Item(99, 12).Value = 1052.8334  # BSM!L99: 1053 -> 1052.8334
$ws.Cells.Item(99, 14).Value = -4048.8334  # BSM!N99: -4049 -> -4048.8334

$ws.Cells.Item(105, 8).Value = 2491.6  # BSM!H105: 2338.739 -> 2491.6
$ws.Cells.Item(105, 9).Value = 2455.2856  # BSM!I105: 2166.5264 -> 2455.2856
$ws.Cells.Item(105, 10).Value = 3000  # BSM!J105: 3156.75 -> 3000
$ws.Cells.Item(105, 11).Value = 2455.2856  # BSM!K105: 2166.5264 -> 2455.2856
$ws.Cells.Item(105, 12).Value = 3000  # BSM!L105: 3156.75 -> 3000
$ws.Cells.Item(105, 13).Value = -708.2856000000002  # BSM!M105: -419.5264000000002 -> -708.2856000000002
$ws.Cells.Item(105, 14).Value = -6494  # BSM!N105: -6650.75 -> -6494

$ws.Cells.Item(132, 8).Value = 50000  # BSM!H132: 0 -> 50000
$ws.Cells.Item(132, 10).Value = 50000  # BSM!J132: 0 -> 50000
$ws.Cells.Item(132, 12).Value = 50000  # BSM!L132: 0 -> 50000
$ws.Cells.Item(132, 14).Value = -60120  # BSM!N132: None -> -60120

$ws.Cells.Item(134, 8).Value = 5175.185  # BSM!H134: 5557.2 -> 5175.185
$ws.Cells.Item(134, 9).Value = 6446.85  # BSM!I134: 6767.737 -> 6446.85
$ws.Cells.Item(134, 10).Value = 1541.8572  # BSM!J134: 1723.8334 -> 1541.8572
$ws.Cells.Item(134, 11).Value = 19340.55  # BSM!K134: 20303.211 -> 19340.55
$ws.Cells.Item(134, 12).Value = 4625.571599999999  # BSM!L134: 5171.5002 -> 4625.571599999999
$ws.Cells.Item(134, 13).Value = -16805.55  # BSM!M134: -17768.211 -> -16805.55
$ws.Cells.Item(134, 14).Value = -9695.571599999999  # BSM!N134: -10241.5002 -> -9695.571599999999

$ws.Cells.Item(137, 8).Value = 61714.145  # BSM!H137: 61860 -> 61714.145
$ws.Cells.Item(137, 10).Value = 61714.145  # BSM!J137: 61860 -> 61714.145
$ws.Cells.Item(137, 12).Value = 61714.145  # BSM!L137: 61860 -> 61714.145
$ws.Cells.Item(137, 14).Value = -71914.14499999999  # BSM!N137: -72060 -> -71914.14499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 5610.909  # CRP!H122: 6532.4 -> 5610.909
$ws.Cells.Item(122, 9).Value = 4213.625  # CRP!I122: 4549.6665 -> 4213.625
$ws.Cells.Item(122, 10).Value = 9337  # CRP!J122: 9506.5 -> 9337
$ws.Cells.Item(122, 11).Value = 12640.875  # CRP!K122: 13648.9995 -> 12640.875
$ws.Cells.Item(122, 12).Value = 28011  # CRP!L122: 28519.5 -> 28011
$ws.Cells.Item(122, 13).Value = -10190.875  # CRP!M122: -11198.9995 -> -10190.875
$ws.Cells.Item(122, 14).Value = -32911  # CRP!N122: -33419.5 -> -32911

$ws.Cells.Item(134, 8).Value = 1596.5883  # CRP!H134: 1803.4615 -> 1596.5883
$ws.Cells.Item(134, 9).Value = 1549.4667  # CRP!I134: 1703.75 -> 1549.4667
$ws.Cells.Item(134, 10).Value = 1950  # CRP!J134: 3000 -> 1950
$ws.Cells.Item(134, 11).Value = 4648.4001  # CRP!K134: 5111.25 -> 4648.4001
$ws.Cells.Item(134, 12).Value = 5850  # CRP!L134: 9000 -> 5850
$ws.Cells.Item(134, 13).Value = -2113.4001  # CRP!M134: -2576.25 -> -2113.4001
$ws.Cells.Item(134, 14).Value = -10920  # CRP!N134: -14070 -> -10920

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1250249.9  # CUL!H4: 1428771.2 -> 1250249.9

$ws.Cells.Item(32, 8).Value = 2999  # CUL!H32: 1499.6666 -> 2999
$ws.Cells.Item(32, 10).Value = 2999  # CUL!J32: 1499.6666 -> 2999
$ws.Cells.Item(32, 12).Value = 8997  # CUL!L32: 4498.9998 -> 8997
$ws.Cells.Item(32, 14).Value = -9563  # CUL!N32: -5064.9998 -> -9563

$ws.Cells.Item(131, 8).Value = 16644.436  # CUL!H131: 18560.902 -> 16644.436
$ws.Cells.Item(131, 10).Value = 17366.227  # CUL!J131: 19473.514 -> 17366.227
$ws.Cells.Item(131, 12).Value = 52098.681  # CUL!L131: 58420.542 -> 52098.681
$ws.Cells.Item(131, 14).Value = -62178.681  # CUL!N131: -68500.542 -> -62178.681

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 3810.8572  # GSM!H70: 3883 -> 3810.8572
$ws.Cells.Item(70, 10).Value = 4126  # GSM!J70: 4500 -> 4126
$ws.Cells.Item(70, 12).Value = 4126  # GSM!L70: 4500 -> 4126
$ws.Cells.Item(70, 14).Value = -4666  # GSM!N70: -5040 -> -4666

$ws.Cells.Item(73, 8).Value = 3810.8572  # GSM!H73: 3883 -> 3810.8572
$ws.Cells.Item(73, 10).Value = 4126  # GSM!J73: 4500 -> 4126
$ws.Cells.Item(73, 12).Value = 4126  # GSM!L73: 4500 -> 4126
$ws.Cells.Item(73, 14).Value = -5998  # GSM!N73: -6372 -> -5998

$ws.Cells.Item(102, 8).Value = 2358.647  # GSM!H102: 2225.7222 -> 2358.647
$ws.Cells.Item(102, 9).Value = 2296.6428  # GSM!I102: 2066.625 -> 2296.6428
$ws.Cells.Item(102, 10).Value = 2648  # GSM!J102: 3498.5 -> 2648
$ws.Cells.Item(102, 11).Value = 2296.6428  # GSM!K102: 2066.625 -> 2296.6428
$ws.Cells.Item(102, 12).Value = 2648  # GSM!L102: 3498.5 -> 2648
$ws.Cells.Item(102, 13).Value = -674.6428000000001  # GSM!M102: -444.625 -> -674.6428000000001
$ws.Cells.Item(102, 14).Value = -5892  # GSM!N102: -6742.5 -> -5892

$ws.Cells.Item(123, 8).Value = 18106  # GSM!H123: 13496.818 -> 18106
$ws.Cells.Item(123, 10).Value = 18106  # GSM!J123: 13496.818 -> 18106
$ws.Cells.Item(123, 12).Value = 18106  # GSM!L123: 13496.818 -> 18106
$ws.Cells.Item(123, 14).Value = -23006  # GSM!N123: -18396.818 -> -23006

$ws.Cells.Item(126, 8).Value = 1887814.8  # GSM!H126: 2177867 -> 1887814.8
$ws.Cells.Item(126, 9).Value = 6947508.5  # GSM!I126: 11114733 -> 6947508.5
$ws.Cells.Item(126, 10).Value = 47926.137  # GSM!J126: 50041.668 -> 47926.137
$ws.Cells.Item(126, 11).Value = 20842525.5  # GSM!K126: 33344199 -> 20842525.5
$ws.Cells.Item(126, 12).Value = 143778.411  # GSM!L126: 150125.004 -> 143778.411
$ws.Cells.Item(126, 13).Value = -20840055.5  # GSM!M126: -33341729 -> -20840055.5
$ws.Cells.Item(126, 14).Value = -148718.411  # GSM!N126: -155065.004 -> -148718.411

$ws.Cells.Item(132, 8).Value = 2749967.2  # GSM!H132: 2264652.5 -> 2749967.2
$ws.Cells.Item(132, 9).Value = 3848524.5  # GSM!I132: 2960368.8 -> 3848524.5
$ws.Cells.Item(132, 11).Value = 11545573.5  # GSM!K132: 8881106.399999999 -> 11545573.5
$ws.Cells.Item(132, 13).Value = -11543043.5  # GSM!M132: -8878576.399999999 -> -11543043.5

$ws.Cells.Item(138, 8).Value = 53029  # GSM!H138: 52929 -> 53029
$ws.Cells.Item(138, 10).Value = 53029  # GSM!J138: 52929 -> 53029
$ws.Cells.Item(138, 12).Value = 53029  # GSM!L138: 52929 -> 53029
$ws.Cells.Item(138, 14).Value = -63309  # GSM!N138: -63209 -> -63309

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(36, 8).Value = 40000  # LTW!H36: 0 -> 40000
$ws.Cells.Item(36, 10).Value = 40000  # LTW!J36: 0 -> 40000
$ws.Cells.Item(36, 12).Value = 40000  # LTW!L36: 0 -> 40000
$ws.Cells.Item(36, 14).Value = -41124  # LTW!N36: None -> -41124

$ws.Cells.Item(46, 8).Value = 1418.0625  # LTW!H46: 1552.7333 -> 1418.0625
$ws.Cells.Item(46, 9).Value = 706.6667  # LTW!I46: 782.5 -> 706.6667
$ws.Cells.Item(46, 10).Value = 2332.7144  # LTW!J46: 2433 -> 2332.7144
$ws.Cells.Item(46, 11).Value = 706.6667  # LTW!K46: 782.5 -> 706.6667
$ws.Cells.Item(46, 12).Value = 2332.7144  # LTW!L46: 2433 -> 2332.7144
$ws.Cells.Item(46, 13).Value = -518.6667  # LTW!M46: -594.5 -> -518.6667
$ws.Cells.Item(46, 14).Value = -2708.7144  # LTW!N46: -2809 -> -2708.7144

$ws.Cells.Item(82, 8).Value = 1887.2307  # LTW!H82: 1969.5834 -> 1887.2307
$ws.Cells.Item(82, 9).Value = 1306.4  # LTW!I82: 1351.6666 -> 1306.4
$ws.Cells.Item(82, 11).Value = 1306.4  # LTW!K82: 1351.6666 -> 1306.4
$ws.Cells.Item(82, 13).Value = -945.4000000000001  # LTW!M82: -990.6666 -> -945.4000000000001

$ws.Cells.Item(85, 8).Value = 1887.2307  # LTW!H85: 1969.5834 -> 1887.2307
$ws.Cells.Item(85, 9).Value = 1306.4  # LTW!I85: 1351.6666 -> 1306.4
$ws.Cells.Item(85, 11).Value = 1306.4  # LTW!K85: 1351.6666 -> 1306.4
$ws.Cells.Item(85, 13).Value = -58.40000000000009  # LTW!M85: -103.6666 -> -58.40000000000009

$ws.Cells.Item(122, 8).Value = 6836.9375  # LTW!H122: 7421.2856 -> 6836.9375
$ws.Cells.Item(122, 9).Value = 7043.6665  # LTW!I122: 6877.5557 -> 7043.6665
$ws.Cells.Item(122, 10).Value = 6571.143  # LTW!J122: 8400 -> 6571.143
$ws.Cells.Item(122, 11).Value = 21130.9995  # LTW!K122: 20632.6671 -> 21130.9995
$ws.Cells.Item(122, 12).Value = 19713.429  # LTW!L122: 25200 -> 19713.429
$ws.Cells.Item(122, 13).Value = -18680.9995  # LTW!M122: -18182.6671 -> -18680.9995
$ws.Cells.Item(122, 14).Value = -24613.429  # LTW!N122: -30100 -> -24613.429

$ws.Cells.Item(132, 8).Value = 2099.0344  # LTW!H132: 2072.7693 -> 2099.0344
$ws.Cells.Item(132, 9).Value = 1257.5714  # LTW!I132: 1300.5 -> 1257.5714
$ws.Cells.Item(132, 10).Value = 2366.7727  # LTW!J132: 2304.45 -> 2366.7727
$ws.Cells.Item(132, 11).Value = 3772.7142  # LTW!K132: 3901.5 -> 3772.7142
$ws.Cells.Item(132, 12).Value = 7100.3181  # LTW!L132: 6913.349999999999 -> 7100.3181
$ws.Cells.Item(132, 13).Value = -1242.7142  # LTW!M132: -1371.5 -> -1242.7142
$ws.Cells.Item(132, 14).Value = -12160.3181  # LTW!N132: -11973.35 -> -12160.3181

$ws.Cells.Item(136, 8).Value = 4756.174  # LTW!H136: 5120.1 -> 4756.174
$ws.Cells.Item(136, 9).Value = 3867.1667  # LTW!I136: 4101.1875 -> 3867.1667
$ws.Cells.Item(136, 10).Value = 7956.6  # LTW!J136: 9195.75 -> 7956.6
$ws.Cells.Item(136, 11).Value = 11601.5001  # LTW!K136: 12303.5625 -> 11601.5001
$ws.Cells.Item(136, 12).Value = 23869.8  # LTW!L136: 27587.25 -> 23869.8
$ws.Cells.Item(136, 13).Value = -9051.500100000001  # LTW!M136: -9753.5625 -> -9051.500100000001
$ws.Cells.Item(136, 14).Value = -28969.8  # LTW!N136: -32687.25 -> -28969.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1150.75  # WVR!H100: 1152.75 -> 1150.75
$ws.Cells.Item(100, 9).Value = 934.3333  # WVR!I100: 937 -> 934.3333
$ws.Cells.Item(100, 11).Value = 1868.6666  # WVR!K100: 1874 -> 1868.6666
$ws.Cells.Item(100, 13).Value = -1327.6666  # WVR!M100: -1333 -> -1327.6666

$ws.Cells.Item(132, 8).Value = 1608.4642  # WVR!H132: 1963.3334 -> 1608.4642
$ws.Cells.Item(132, 9).Value = 1127.875  # WVR!I132: 1368.3529 -> 1127.875
$ws.Cells.Item(132, 11).Value = 3383.625  # WVR!K132: 4105.0587 -> 3383.625
$ws.Cells.Item(132, 13).Value = -853.625  # WVR!M132: -1575.0587 -> -853.625

$ws.Cells.Item(135, 8).Value = 139857.2  # WVR!H135: 136547.5 -> 139857.2
$ws.Cells.Item(135, 10).Value = 139857.2  # WVR!J135: 136547.5 -> 139857.2
$ws.Cells.Item(135, 12).Value = 139857.2  # WVR!L135: 136547.5 -> 139857.2
$ws.Cells.Item(135, 14).Value = -149997.2  # WVR!N135: -146687.5 -> -149997.2

$ws.Cells.Item(136, 8).Value = 19159324  # WVR!H136: 17363274 -> 19159324
$ws.Cells.Item(136, 9).Value = 34724652  # WVR!I136: 32682050 -> 34724652
$ws.Cells.Item(136, 10).Value = 1999.4615  # WVR!J136: 1991.9333 -> 1999.4615
$ws.Cells.Item(136, 11).Value = 104173956  # WVR!K136: 98046150 -> 104173956
$ws.Cells.Item(136, 12).Value = 5998.3845  # WVR!L136: 5975.7999 -> 5998.3845
$ws.Cells.Item(136, 13).Value = -104171406  # WVR!M136: -98043600 -> -104171406
$ws.Cells.Item(136, 14).Value = -11098.3845  # WVR!N136: -11075.7999 -> -11098.3845

$ws.Cells.Item(138, 8).Value = 0  # WVR!H138: 62000 -> 0
$ws.Cells.Item(138, 10).Value = 0  # WVR!J138: 62000 -> 0
$ws.Cells.Item(138, 12).Value = 0  # WVR!L138: 62000 -> 0
$ws.Cells.Item(138, 14).ClearContents()  # WVR!N138: -72280 -> (removed)
